$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G15"  = "2025/2026, 2023/2024"
    "G19"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G21"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G22"  = "2025/2026, Eman_mohamed@med.asu.edu.eg"
    "G37"  = "2025/2026, 2023/2024"
    "G41"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G43"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G44"  = "2025/2026, Eman_mohamed@med.asu.edu.eg"
    "G60"  = "2025/2026, 2026/2027"
    "G63"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G64"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G65"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G66"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G82"  = "2025/2026, 2026/2027"
    "G85"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G86"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G87"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G88"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G89"  = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G104" = "2025/2026, 2026/2027"
    "G106" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G107" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G110" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G111" = "marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026, yassmen.ahmed@med.asu.edu.eg"
    "G126" = "2025/2026, 2026/2027"
    "G128" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G129" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G132" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G133" = "marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026, yassmen.ahmed@med.asu.edu.eg"
    "G150" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G153" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G155" = "marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026, yassmen.ahmed@med.asu.edu.eg"
    "G172" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G175" = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    "G177" = "marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, 2025/2026, yassmen.ahmed@med.asu.edu.eg"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
